$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name swap: Marruecos now listed before Uzbekistan in the shared-string table.
# The underlying data row (A63) that used to be "Uzbekistan" becomes "Marruecos" (with fresh stats),
# and the row below (A64) becomes "Uzbekistan" (carrying the old Uzbekistan stats down).
$ws.Range("A63").Value = "Marruecos"
$ws.Range("A64").Value = "Uzbekistan"

# --- Updated statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
# Row 4
$ws.Range("B4").Value = 4832187
$ws.Range("C4").Value = 18540
$ws.Range("D4").Value = 2389698
$ws.Range("E4").Value = 2283927
$ws.Range("G4").Value = 197
$ws.Range("H4").Value = 158562

# Row 5
$ws.Range("B5").Value = 2736298
$ws.Range("C5").Value = 2621
$ws.Range("E5").Value = 758021
$ws.Range("G5").Value = 96
$ws.Range("H5").Value = 94226

# Row 6
$ws.Range("B6").Value = 1852156
$ws.Range("C6").Value = 47454
$ws.Range("D6").Value = 1229171
$ws.Range("E6").Value = 584016
$ws.Range("G6").Value = 808
$ws.Range("H6").Value = 38969

# Row 21
$ws.Range("B21").Value = 212060
$ws.Range("C21").Value = 598
$ws.Range("E21").Value = 9234

# Row 36
$ws.Range("B36").Value = 74102
$ws.Range("C36").Value = 1287
$ws.Range("D36").Value = 47551
$ws.Range("E36").Value = 26005
$ws.Range("G36").Value = 10
$ws.Range("H36").Value = 546

# Row 62
$ws.Range("B62").Value = 26208
$ws.Range("C62").Value = 46
$ws.Range("E62").Value = 1081

# Row 63
$ws.Range("B63").Value = 26196
$ws.Range("C63").Value = 659
$ws.Range("D63").Value = 18968
$ws.Range("E63").Value = 6827
$ws.Range("G63").Value = 19
$ws.Range("H63").Value = 401

# Row 64
$ws.Range("B64").Value = 26066
$ws.Range("C64").Value = 730
$ws.Range("D64").Value = 16838
$ws.Range("E64").Value = 9071
$ws.Range("G64").Value = 6
$ws.Range("H64").Value = 157

# Row 95
$ws.Range("D95").Value = 5498
$ws.Range("E95").Value = 1248

# Row 102
$ws.Range("D102").Value = 1837
$ws.Range("E102").Value = 3160

# Row 131
$ws.Range("B131").Value = 1973
$ws.Range("C131").Value = 27
$ws.Range("E131").Value = 1305
$ws.Range("G131").Value = 1
$ws.Range("H131").Value = 14

# Row 136
$ws.Range("D136").Value = 863
$ws.Range("E136").Value = 372
$ws.Range("G136").Value = 2
$ws.Range("H136").Value = 499

# Row 142
$ws.Range("B142").Value = 1214
$ws.Range("C142").Value = 7
$ws.Range("D142").Value = 696
$ws.Range("E142").Value = 440
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 78

# Row 162
$ws.Range("D162").Value = 374
$ws.Range("E162").Value = 262

# --- Last-updated timestamp
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 20:09"
